$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Motif" shifts from C to D, etc.)
$ws.Columns("C:C").Insert()

# Header for the new "Factor" column
$ws.Range("C1").Value = "Factor"

# Populate the Factor values (gene/TF symbol parsed from the Motif name)
# for the rows that already carried per-motif metadata (rows 2-14).
$ws.Range("C2").Value = "IRF1"
$ws.Range("C3").Value = "Z354A"
$ws.Range("C4").Value = "MAZ"
$ws.Range("C5").Value = "THAP1"
$ws.Range("C6").Value = "EGR1"
$ws.Range("C7").Value = "SP1"
$ws.Range("C8").Value = "PATZ1"
$ws.Range("C9").Value = "VEZF1"
$ws.Range("C10").Value = "KLF6"
$ws.Range("C11").Value = "KLF1"
$ws.Range("C12").Value = "RREB1"
$ws.Range("C13").Value = "MGA"
$ws.Range("C14").Value = "Z324A"

# Refresh the cached sort-state metadata so it reflects the new column
# layout (data stays in the same order - this only rewrites the stored
# sort range/condition references).
$sortRange = $ws.Range("A2:H134")
$keyRange = $ws.Range("G2:G134")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Restore the active selection to the new, empty Factor cell for row 15.
$ws.Range("C15").Select()
